$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'copper compression shorts',
    'copper compression shorts women',
    'copper compression tights women',
    'copper compression womens shorts',
    'copper compressions',
    'copper fit clothes for women',
    'copper fit compression leggings',
    'copper fit compression pants',
    'copper fit hamstring',
    'copper fit hip compression',
    'copper fit leggings for women',
    'copper fit leggings women',
    'copper fit lower back',
    'copper fit pants',
    'copper fit pants for women',
    'copper fit pants women',
    'copper fit tights for women',
    'copper fit underwear',
    'copper fitness',
    'copper infused clothing',
    'copper key brand leggings',
    'copper key clothing',
    'copper key underwear',
    'copper leggings',
    'copper leggings compression women',
    'copper leggings for women',
    'copper pants',
    'copper shorts for women',
    'copper tights',
    'copper tops for women',
    'copper underwear',
    'copper wear compression',
    'copper yoga pants',
    'core compression band',
    'core daily 1 women 50',
    'core de force apparel',
    'core leggings',
    'core power shakes',
    'core secrets ball',
    'core10 leggings',
    'cothes basket',
    'cotton active running bike leggings-athletic exercise yoga walking 7/3 shorts',
    'country apperal for woman',
    'country leggings',
    'couple workout gear',
    'cow apparel women',
    'cow apperal for women',
    'cow balls',
    'cow body wash',
    'cow hip lift',
    'cow pants women',
    'cow shorts',
    'cow shorts for women',
    'cow thigh high',
    'cow tights',
    'cow web',
    'cozy bluefish costume',
    'cracked skin leggings',
    'cracked tights',
    'crazy leggings for women',
    'crazy legs',
    'crazy sucks',
    'crazy tights for women',
    'crazy yoga leggings',
    'crazy yoga loose running shorts',
    'crazy yoga shorts',
    'crew 3 in 1 body wash',
    'crop leggings for women',
    'crop leggings workout',
    'crop pants for men',
    'crop tights',
    'crop tights for women',
    'crop workout pants',
    'cropped leggings white',
    'cross band leggings',
    'cross check',
    'cross country pants women',
    'cross country running gear',
    'cross country skiing apparel',
    'cross country skiing clothes women',
    'cross country skiing pants',
    'cross country skiing pants men',
    'cross dress men to women clothes',
    'cross fit pants women',
    'cross knee brace',
    'cross leggings for women',
    'crossfit apparel for women',
    'crossfit back support',
    'crossfit compression pants men',
    'crossfit games 2018 womens apparel',
    'crossfit knee support',
    'crossfit leggings',
    'crossfit leggings women',
    'crossfit pants',
    'crossfit pants for women',
    'crossfit recovery',
    'crossfit women',
    'crx leggings',
    'crx yoga',
    'travel exercise'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
